$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rewrite the header row (row 1) into its new layout: "Nam hoc/Khoi/Lop"
# move up near the front (after STT), "Ma so hoc sinh" and "Trang thai"
# columns are dropped, and "Fist Name"/"Last Name" become "Ho"/"Ten".
$ws.Range("A1").Value = 'STT'
$ws.Range("B1").Value = 'Năm học'
$ws.Range("C1").Value = 'Khối'
$ws.Range("D1").Value = 'Lớp'
$ws.Range("E1").Value = 'Họ'
$ws.Range("F1").Value = 'Tên'
$ws.Range("G1").Value = 'Năm sinh'
$ws.Range("H1").Value = 'Giới tính'
$ws.Range("I1").Value = 'Dân tộc'
$ws.Range("J1").Value = 'Ngày vào trường'
$ws.Range("K1").Value = 'Số điện thoại'
$ws.Range("L1").Value = 'Địa chỉ'
$ws.Range("M1").Value = 'Cha'
$ws.Range("N1").Value = 'Mẹ'
$ws.Range("O1").Value = 'Quan hệ khác'
$ws.Range("P1").Value = 'Họ tên cha'
$ws.Range("Q1").Value = 'Năm sinh cha'
$ws.Range("R1").Value = 'Số điện thoại cha'
$ws.Range("S1").Value = 'Nghề nghiệp cha'
$ws.Range("T1").Value = 'Họ tên mẹ'
$ws.Range("U1").Value = 'Năm sinh mẹ'
$ws.Range("V1").Value = 'Số điện thoại mẹ'
$ws.Range("W1").Value = 'Nghề nghiệp mẹ'
$ws.Range("X1").Value = 'Họ tên quan hệ khác'
$ws.Range("Y1").Value = 'Năm sinh quan hệ khác'
$ws.Range("Z1").Value = 'Số điện thoại quan hệ khác'
$ws.Range("AA1").Value = 'Nghề nghiệp quan hệ khác'

# The old sheet had two extra trailing columns (AB, AC) that are no longer
# used now that the row only spans through AA - clear them out so the used
# range / dimension shrinks back down to A1:AA1.
$ws.Range("AB1:AC1").ClearContents()

# Update the remembered selection to match the saved view.
$ws.Range("K10").Select()
